$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.33083176612854
$ws.Range("B1").Value = 4.3038010597229
$ws.Range("C1").Value = 3.271773338317871
$ws.Range("D1").Value = 1.021642446517944
$ws.Range("E1").Value = 0.4989390671253204
